$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117, shifting existing rows 117-118 down to 118-119
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new weekly data point
$ws.Cells.Item(117, 1).Value = 5
$ws.Cells.Item(117, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(117, 3).Value = "Maule"
$ws.Cells.Item(117, 4).Value = 44890
$ws.Cells.Item(117, 5).Value = 7
$ws.Cells.Item(117, 6).Value = 100112022
$ws.Cells.Item(117, 7).Value = "Arveja Verde"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 300
$ws.Cells.Item(117, 11).Value = 23000
$ws.Cells.Item(117, 12).Value = 23000
$ws.Cells.Item(117, 13).Value = 23000
$ws.Cells.Item(117, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(117, 15).Value = "Región del Maule"
$ws.Cells.Item(117, 16).Value = 920
$ws.Cells.Item(117, 17).Value = 25
$ws.Cells.Item(117, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D (copy from the row below)
$ws.Cells.Item(118, 4).Copy()
$ws.Cells.Item(117, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(117, 4).Value = 44890

# The row that was previously at 117 is now at 118; its Fecha changes from 44496 to 44263
$ws.Cells.Item(118, 4).Value = 44263
